$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 47898.523
$ws.Range("I33").Value = 55717.832
$ws.Range("J33").Value = 982.6667
$ws.Range("K33").Value = 55717.832
$ws.Range("L33").Value = 982.6667
$ws.Range("M33").Value = -55488.832
$ws.Range("N33").Value = -1440.6667

$ws.Range("H41").Value = 3996.6667
$ws.Range("I41").Value = 3996.6667
$ws.Range("K41").Value = 3996.6667
$ws.Range("M41").Value = -3556.6667

$ws.Range("H43").Value = 6316
$ws.Range("I43").Value = 6246.8335
$ws.Range("K43").Value = 6246.8335
$ws.Range("M43").Value = -6177.8335

$ws.Range("H70").Value = 2332.5454
$ws.Range("J70").Value = 2558
$ws.Range("L70").Value = 7674
$ws.Range("N70").Value = -8214

$ws.Range("H73").Value = 2332.5454
$ws.Range("J73").Value = 2558
$ws.Range("L73").Value = 7674
$ws.Range("N73").Value = -9546

$ws.Range("H76").Value = 4966.5835
$ws.Range("I76").Value = 4759.9
$ws.Range("K76").Value = 4759.9
$ws.Range("M76").Value = -4444.9

$ws.Range("H79").Value = 4966.5835
$ws.Range("I79").Value = 4759.9
$ws.Range("K79").Value = 4759.9
$ws.Range("M79").Value = -3667.9

$ws.Range("H80").Value = 417.1
$ws.Range("I80").Value = 183.5
$ws.Range("K80").Value = 550.5
$ws.Range("M80").Value = 447.5

$ws.Range("H83").Value = 417.1
$ws.Range("I83").Value = 183.5
$ws.Range("K83").Value = 1651.5
$ws.Range("M83").Value = 3340.5

$ws.Range("H86").Value = 5114.5713
$ws.Range("I86").Value = 2724.5
$ws.Range("J86").Value = 8301.333000000001
$ws.Range("K86").Value = 2724.5
$ws.Range("L86").Value = 8301.333000000001
$ws.Range("M86").Value = -1601.5
$ws.Range("N86").Value = -10547.333

$ws.Range("H89").Value = 5114.5713
$ws.Range("I89").Value = 2724.5
$ws.Range("J89").Value = 8301.333000000001
$ws.Range("K89").Value = 13622.5
$ws.Range("L89").Value = 41506.665
$ws.Range("M89").Value = -8006.5
$ws.Range("N89").Value = -52738.665

$ws.Range("H125").Value = 1113.5555
$ws.Range("I125").Value = 999
$ws.Range("J125").Value = 1127.875
$ws.Range("K125").Value = 8991
$ws.Range("L125").Value = 10150.875
$ws.Range("M125").Value = -6531
$ws.Range("N125").Value = -15070.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7594.516
$ws.Range("I32").Value = 7797.6665
$ws.Range("K32").Value = 7797.6665
$ws.Range("M32").Value = -7510.6665

$ws.Range("H97").Value = 2818.074
$ws.Range("I97").Value = 2947.52
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 2947.52
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -2451.52
$ws.Range("N97").Value = -2192

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H114").Value = 29999.666
$ws.Range("J114").Value = 29999.666
$ws.Range("L114").Value = 29999.666
$ws.Range("N114").Value = -38677.666

$ws.Range("H122").Value = 2629.2856
$ws.Range("I122").Value = 2651
$ws.Range("K122").Value = 7953
$ws.Range("M122").Value = -5503

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5238.75
$ws.Range("I20").Value = 7539.0713
$ws.Range("K20").Value = 7539.0713
$ws.Range("M20").Value = -7292.0713

$ws.Range("H86").Value = 5262.7
$ws.Range("I86").Value = 880.9375
$ws.Range("K86").Value = 880.9375
$ws.Range("M86").Value = 242.0625

$ws.Range("H89").Value = 5262.7
$ws.Range("I89").Value = 880.9375
$ws.Range("K89").Value = 4404.6875
$ws.Range("M89").Value = 1211.3125

$ws.Range("H134").Value = 2767.5557
$ws.Range("I134").Value = 2636.2354
$ws.Range("K134").Value = 7908.706200000001
$ws.Range("M134").Value = -5373.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1599
$ws.Range("I58").Value = 1599
$ws.Range("K58").Value = 1599
$ws.Range("M58").Value = -1396

$ws.Range("H86").Value = 3250.6667
$ws.Range("I86").Value = 2637.6667
$ws.Range("J86").Value = 3863.6667
$ws.Range("K86").Value = 2637.6667
$ws.Range("L86").Value = 3863.6667
$ws.Range("M86").Value = -1514.6667
$ws.Range("N86").Value = -6109.6667

$ws.Range("H89").Value = 3250.6667
$ws.Range("I89").Value = 2637.6667
$ws.Range("J89").Value = 3863.6667
$ws.Range("K89").Value = 13188.3335
$ws.Range("L89").Value = 19318.3335
$ws.Range("M89").Value = -7572.333500000001
$ws.Range("N89").Value = -30550.3335

$ws.Range("H132").Value = 1373.6364
$ws.Range("I132").Value = 1361.9
$ws.Range("J132").Value = 1491
$ws.Range("K132").Value = 4085.7
$ws.Range("L132").Value = 4473
$ws.Range("M132").Value = -1555.7
$ws.Range("N132").Value = -9533

$ws.Range("H136").Value = 1599
$ws.Range("I136").Value = 1599
$ws.Range("K136").Value = 4797
$ws.Range("M136").Value = -2247

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3792.6365
$ws.Range("J80").Value = 4361.2856
$ws.Range("L80").Value = 13083.8568
$ws.Range("N80").Value = -14955.8568

$ws.Range("H83").Value = 3792.6365
$ws.Range("J83").Value = 4361.2856
$ws.Range("L83").Value = 39251.5704
$ws.Range("N83").Value = -48611.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6738.2617
$ws.Range("I70").Value = 5373.3335
$ws.Range("J70").Value = 6843.2563
$ws.Range("K70").Value = 5373.3335
$ws.Range("L70").Value = 6843.2563
$ws.Range("M70").Value = -5103.3335
$ws.Range("N70").Value = -7383.2563

$ws.Range("H73").Value = 6738.2617
$ws.Range("I73").Value = 5373.3335
$ws.Range("J73").Value = 6843.2563
$ws.Range("K73").Value = 5373.3335
$ws.Range("L73").Value = 6843.2563
$ws.Range("M73").Value = -4437.3335
$ws.Range("N73").Value = -8715.256300000001

$ws.Range("H97").Value = 748.75
$ws.Range("I97").Value = 499
$ws.Range("J97").Value = 1165
$ws.Range("K97").Value = 499
$ws.Range("L97").Value = 1165
$ws.Range("M97").Value = -3
$ws.Range("N97").Value = -2157

$ws.Range("H113").Value = 3417.7
$ws.Range("I113").Value = 3359.625
$ws.Range("K113").Value = 3359.625
$ws.Range("M113").Value = -1189.625

$ws.Range("H126").Value = 14166.667
$ws.Range("I126").Value = 28000
$ws.Range("J126").Value = 12437.5
$ws.Range("K126").Value = 84000
$ws.Range("L126").Value = 37312.5
$ws.Range("M126").Value = -81530
$ws.Range("N126").Value = -42252.5

$ws.Range("H132").Value = 1956.2
$ws.Range("I132").Value = 1773.3846
$ws.Range("K132").Value = 5320.1538
$ws.Range("M132").Value = -2790.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2523.0625
$ws.Range("I68").Value = 2523.0625
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2523.0625
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 2523.0625
$ws.Range("I71").Value = 2523.0625
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 12615.3125
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H93").Value = 26032.875
$ws.Range("J93").Value = 200495
$ws.Range("L93").Value = 200495
$ws.Range("N93").Value = -202991

$ws.Range("H132").Value = 5126.909
$ws.Range("I132").Value = 3343.7144
$ws.Range("J132").Value = 8247.5
$ws.Range("K132").Value = 10031.1432
$ws.Range("L132").Value = 24742.5
$ws.Range("M132").Value = -7501.143199999999
$ws.Range("N132").Value = -29802.5

$ws.Range("H136").Value = 31253628
$ws.Range("I136").Value = 3245.6897
$ws.Range("J136").Value = 333340670
$ws.Range("K136").Value = 9737.069100000001
$ws.Range("L136").Value = 1000022010
$ws.Range("M136").Value = -7187.069100000001
$ws.Range("N136").Value = -1000027110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 14448.75
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 14448.75
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184

$ws.Range("H96").Value = 2295.6365
$ws.Range("I96").Value = 2417
$ws.Range("J96").Value = 2150
$ws.Range("K96").Value = 2417
$ws.Range("L96").Value = 2150
$ws.Range("M96").Value = -1044
$ws.Range("N96").Value = -4896

$ws.Range("H100").Value = 1236.75
$ws.Range("I100").Value = 818.6
$ws.Range("K100").Value = 1637.2
$ws.Range("M100").Value = -1096.2

$ws.Range("H107").Value = 1446.5
$ws.Range("I107").Value = 1208
$ws.Range("J107").Value = 1804.25
$ws.Range("K107").Value = 3624
$ws.Range("L107").Value = 5412.75
$ws.Range("M107").Value = -1704
$ws.Range("N107").Value = -9252.75

$ws.Range("H132").Value = 2818.3635
$ws.Range("I132").Value = 2818.3635
$ws.Range("K132").Value = 2818.3635
$ws.Range("M132").Value = -5925.0905
